$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 into I1:J1 (bold/centered/bordered header xf), then set labels
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I2:J62 ("I0" / "IF" columns) with the new data
$data = New-Object "object[,]" 61,2
$data[0,0] = 8
$data[0,1] = 8
$data[1,0] = 6
$data[1,1] = 7
$data[2,0] = 6
$data[2,1] = 6
$data[3,0] = 7
$data[3,1] = 8
$data[4,0] = 7
$data[4,1] = 7
$data[5,0] = 9
$data[5,1] = 9
$data[6,0] = 8
$data[6,1] = 8
$data[7,0] = 8
$data[7,1] = 8
$data[8,0] = 8
$data[8,1] = 8
$data[9,0] = 11
$data[9,1] = 11
$data[10,0] = 5
$data[10,1] = 5
$data[11,0] = 6
$data[11,1] = 7
$data[12,0] = 4
$data[12,1] = 5
$data[13,0] = 6
$data[13,1] = 6
$data[14,0] = 10
$data[14,1] = 10
$data[15,0] = 7
$data[15,1] = 7
$data[16,0] = 6
$data[16,1] = 7
$data[17,0] = 5
$data[17,1] = 5
$data[18,0] = 6
$data[18,1] = 7
$data[19,0] = 5
$data[19,1] = 5
$data[20,0] = 5
$data[20,1] = 6
$data[21,0] = 7
$data[21,1] = 7
$data[22,0] = 5
$data[22,1] = 6
$data[23,0] = 7
$data[23,1] = 7
$data[24,0] = 6
$data[24,1] = 6
$data[25,0] = 5
$data[25,1] = 5
$data[26,0] = 9
$data[26,1] = 9
$data[27,0] = 4
$data[27,1] = 4
$data[28,0] = 9
$data[28,1] = 9
$data[29,0] = 7
$data[29,1] = 7
$data[30,0] = 7
$data[30,1] = 8
$data[31,0] = 7
$data[31,1] = 7
$data[32,0] = 6
$data[32,1] = 6
$data[33,0] = 1
$data[33,1] = 3
$data[34,0] = 5
$data[34,1] = 6
$data[35,0] = 6
$data[35,1] = 7
$data[36,0] = 8
$data[36,1] = 8
$data[37,0] = 4
$data[37,1] = 5
$data[38,0] = 3
$data[38,1] = 5
$data[39,0] = 5
$data[39,1] = 6
$data[40,0] = 6
$data[40,1] = 6
$data[41,0] = 7
$data[41,1] = 7
$data[42,0] = 5
$data[42,1] = 5
$data[43,0] = 1
$data[43,1] = 3
$data[44,0] = 7
$data[44,1] = 8
$data[45,0] = 7
$data[45,1] = 7
$data[46,0] = 6
$data[46,1] = 7
$data[47,0] = 5
$data[47,1] = 5
$data[48,0] = 7
$data[48,1] = 8
$data[49,0] = 8
$data[49,1] = 8
$data[50,0] = 6
$data[50,1] = 6
$data[51,0] = 6
$data[51,1] = 8
$data[52,0] = 9
$data[52,1] = 9
$data[53,0] = 6
$data[53,1] = 6
$data[54,0] = 8
$data[54,1] = 9
$data[55,0] = 8
$data[55,1] = 8
$data[56,0] = 6
$data[56,1] = 8
$data[57,0] = 9
$data[57,1] = 9
$data[58,0] = 8
$data[58,1] = 8
$data[59,0] = 8
$data[59,1] = 8
$data[60,0] = 1
$data[60,1] = 1

$ws.Range("I2:J62").Value = $data

